$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 581, shifting rows 581:622 down to 582:623.
$ws.Rows("581:581").Insert()

# Populate the new row 581 with the new data point.
# Column A holds a date-look-alike string (e.g. "2026/01/08") that must stay
# a plain text value (matching the sheet's existing inline-string cells)
# rather than being auto-converted into a real Excel date. Temporarily force
# the cell to Text format while assigning, then restore the default
# "Normal" style so no stray number-format override is left behind.
$ws.Range("A581").NumberFormat = "@"
$ws.Range("A581").Value = "2026/01/08"
$ws.Range("A581").Style = "Normal"

$ws.Range("B581").Value = "木"
$ws.Range("C581").Value = 17
$ws.Range("D581").Value = 24
